# Fill in the previously-empty metric-description cells (column B) for the
# metric groups that had no description yet: fgw_meta_info, fgw_resource_usage,
# fgw_bandwidth and fgw_http_latency. These cells are the top-left cell of a
# vertically merged range (B8:B11, B12:B16, B24:B28, B37:B41).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B8").Value  = "元数据，用于与其他指标联合使用，减少重复字段"
$ws.Range("B12").Value = "资源(CPU和内存)使用情况"
$ws.Range("B24").Value = "带宽使用情况"
$ws.Range("B37").Value = "7层HTTP请求的延迟情况"
